$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I0 (col I) and IF (col J), matching the style of H1 (bold + border + centered/top alignment)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# I/J values per data row (row 2..71 -> I0, IF)
$values = @(
    @(5,5),
    @(7,7),
    @(9,9),
    @(8,8),
    @(8,8),
    @(9,9),
    @(10,10),
    @(9,9),
    @(8,8),
    @(9,9),
    @(8,8),
    @(7,7),
    @(9,9),
    @(9,9),
    @(8,8),
    @(8,8),
    @(9,9),
    @(8,8),
    @(8,8),
    @(8,8),
    @(8,8),
    @(8,8),
    @(8,8),
    @(8,8),
    @(8,8),
    @(9,9),
    @(8,8),
    @(8,8),
    @(8,8),
    @(7,7),
    @(7,7),
    @(6,6),
    @(5,6),
    @(9,9),
    @(6,6),
    @(6,7),
    @(7,7),
    @(6,7),
    @(6,6),
    @(6,7),
    @(8,8),
    @(7,7),
    @(7,7),
    @(7,7),
    @(6,6),
    @(8,8),
    @(7,7),
    @(7,7),
    @(8,8),
    @(7,7),
    @(8,8),
    @(8,8),
    @(8,8),
    @(7,7),
    @(8,8),
    @(9,9),
    @(9,9),
    @(8,8),
    @(9,9),
    @(8,8),
    @(8,8),
    @(9,9),
    @(7,7),
    @(7,7),
    @(7,7),
    @(7,7),
    @(6,6),
    @(4,5),
    @(7,7),
    @(3,3)
)

for ($k = 0; $k -lt $values.Count; $k++) {
    $row = 2 + $k
    $pair = $values[$k]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
